$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = -7.43169999999999
$ws.Range("C3").Value = -12.6263
$ws.Range("D5").Value = -8.009200000000007
$ws.Range("D6").Value = -8.109500000000001
$ws.Range("D8").Value = -8.894800000000002
$ws.Range("C9").Value = -10.272
$ws.Range("A11").Value = -21.82529999999999
$ws.Range("B11").Value = 5.6375
$ws.Range("A12").Value = -21.3512
$ws.Range("C13").Value = -13.0571
$ws.Range("C14").Value = -13.85219999999999
$ws.Range("A15").Value = -21.49779999999999
$ws.Range("D17").Value = -8.28999999999999
$ws.Range("C19").Value = -12.52920000000001
$ws.Range("C21").Value = -12.41859999999999
$ws.Range("C22").Value = -10.3826
$ws.Range("B23").Value = 8.577700000000004
$ws.Range("C24").Value = -12.35139999999999
$ws.Range("C26").Value = -11.9453
$ws.Range("A27").Value = -21.8936
$ws.Range("D27").Value = -8.100400000000004
$ws.Range("A28").Value = -21.5933
$ws.Range("B28").Value = 6.167999999999999
$ws.Range("A31").Value = -21.43729999999999
$ws.Range("A32").Value = -21.35489999999998
$ws.Range("B32").Value = 5.320900000000003
$ws.Range("D33").Value = -7.926200000000001
$ws.Range("B34").Value = 9.59320000000001
$ws.Range("A36").Value = -21.44769999999999
$ws.Range("B36").Value = 5.073099999999999
$ws.Range("B37").Value = 9.027600000000001
$ws.Range("A38").Value = -19.51379999999999
$ws.Range("C38").Value = -12.1569
$ws.Range("C41").Value = -12.43970000000001
$ws.Range("B42").Value = 10.3201
$ws.Range("A46").Value = -21.90699999999999
$ws.Range("B49").Value = 5.733000000000001
$ws.Range("C52").Value = -10.97549999999999
$ws.Range("A54").Value = -21.98650000000001
$ws.Range("B54").Value = 5.224799999999998
$ws.Range("A55").Value = -21.5503
$ws.Range("D55").Value = -8.0588
$ws.Range("A56").Value = -21.7742
$ws.Range("C56").Value = -12.6801
$ws.Range("D59").Value = -8.029199999999994
$ws.Range("A67").Value = -21.62709999999996
$ws.Range("A69").Value = -21.71489999999997
$ws.Range("D70").Value = -5.987299999999999
$ws.Range("C71").Value = -12.1546
$ws.Range("A72").Value = -21.80829999999998
$ws.Range("C72").Value = -12.3596
$ws.Range("A73").Value = -19.45860000000001
$ws.Range("B78").Value = 9.405099999999994
$ws.Range("C78").Value = -10.7847
$ws.Range("B80").Value = 8.880000000000003
$ws.Range("D80").Value = -8.264699999999999
$ws.Range("A83").Value = -21.52359999999999
$ws.Range("C83").Value = -12.80279999999999
$ws.Range("C85").Value = -13.167
$ws.Range("A86").Value = -21.94180000000001
$ws.Range("C86").Value = -13.1713
$ws.Range("C90").Value = -10.08130000000001
$ws.Range("A91").Value = -20.60629999999999
$ws.Range("A93").Value = -21.44920000000002
$ws.Range("D95").Value = -7.920100000000004
$ws.Range("C96").Value = -9.879399999999999
$ws.Range("B97").Value = 5.938700000000003
$ws.Range("D97").Value = -8.341500000000011
$ws.Range("D98").Value = -7.850699999999997
$ws.Range("A99").Value = -21.8748
$ws.Range("B99").Value = 5.202799999999994
$ws.Range("B100").Value = 4.3803
$ws.Range("B101").Value = 5.144199999999997
$ws.Range("D102").Value = -8.0031
$ws.Range("C103").Value = -13.44459999999999
$ws.Range("A104").Value = -21.56099999999999
$ws.Range("A105").Value = -19.82559999999999
